$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Config values flipped from "Yes" to "No" (B5, B6)
$ws.Range("B5").Value = "No"
$ws.Range("B6").Value = "No"

# Slightly widen columns A:B (24.109375 -> ~24.140625 chars)
$ws.Columns("A:B").ColumnWidth = 23.333333

# Move the active selection from B5 to B7
[void]$ws.Range("B7").Select()
